$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 31
$ws1.Range("F9").Value = 194
$ws1.Range("F11").Value = 49
$ws1.Range("F13").Value = 1480
$ws1.Range("F15").Value = 2742

# Sheet "全部类型" (All Types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 31
$ws4.Range("F10").Value = 194
$ws4.Range("F12").Value = 49
$ws4.Range("F16").Value = 1480
$ws4.Range("F18").Value = 2742
